# Generate Report for Handback
# Updates the handback-status report with refreshed timestamps and status
# for the latest handback/regeneration pass.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# "Latest HO Xliff Generate Date" column (G) for the 0611575f... file
# is recomputed (rows 2 and 5 both reference this same file/date).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-27 20:17:16"
$wsOverview.Range("G5").Value = "2016-08-27 20:17:16"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Status changed from "ht" (human translation) to "mt" (machine translation)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
# Correspond Handoff Datetime (H) refreshed
$wsZhCn.Range("H2").Value = "2016-08-27 20:17:11"
$wsZhCn.Range("H5").Value = "2016-08-27 20:17:11"
# Correspond Handback DateTime (K) refreshed
$wsZhCn.Range("K2").Value = "2016-08-27 20:17:29"
$wsZhCn.Range("K5").Value = "2016-08-27 20:17:29"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Status changed from "ht" to "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
# Correspond Handoff Datetime (H) refreshed
$wsDeDe.Range("H2").Value = "2016-08-27 20:17:16"
$wsDeDe.Range("H5").Value = "2016-08-27 20:17:16"
# Correspond Handback DateTime (K) refreshed
$wsDeDe.Range("K2").Value = "2016-08-27 20:17:35"
$wsDeDe.Range("K5").Value = "2016-08-27 20:17:35"
